$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(6).Insert()
$ws.Range("A6").Value = "TEK69_3"
$ws.Range("A1:A12").SetPhonetic()
$ws.Range("A4").Select()
